$d = $word.ActiveDocument

# --- Text edits in the "Another..." paragraph (paragraph 5) ---

# 1. "I have the perception though" -> "I had the perception growing up though"
$d.Content.Find.Execute(
    "religious.  I have the perception though that Jews", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "religious.  I had the perception growing up though that Jews", 2)

# 2. "chastised. However, when this movie" -> "chastised. When this movie"
$d.Content.Find.Execute(
    "often chastised. However, when this movie came out,", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "often chastised. When this movie came out,", 2)

# 3. "period of time and for a while" -> "period of time, and for a while"
$d.Content.Find.Execute(
    "very long period of time and for a while", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "very long period of time, and for a while", 2)

# 4. "later years of High School. However," -> "later years of High School because of this residual effect. However,"
$d.Content.Find.Execute(
    "later years of High School. However, I am happy", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "later years of High School because of this residual effect. However, I am happy", 2)

# --- Move the "_GoBack" bookmark from the end of paragraph 3 to wrap paragraph 5 ---

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$p5 = $d.Paragraphs.Item(5)
$newRange = $d.Range($p5.Range.Start, $d.Content.End)
$newRange.Bookmarks.Add("_GoBack")
